$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad" / Changed date) is bumped by one day (45181 -> 45182)
# for every data row (rows 2 through 171).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2()
    if ($current -eq 45181) {
        $cell.Value = 45182
    }
}
